$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Delete the old rows 7-14 (data shrinks from 14 rows to 6 rows total) ---
$ws.Rows("7:14").Delete()

# --- Clear all existing hyperlinks on the sheet; they'll be re-added below for the surviving rows ---
$ws.Cells.Hyperlinks.Delete()

# --- Update data rows 2-6 with the new scraped listings ---
$ws.Range("A2").Value = "2026-01-08 06:31:27"
$ws.Range("B2").Value = "【法人歓迎】プローバステージ制御ソフト開発の見積依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5467295"
$ws.Range("G2").Value = 75
$ws.Range("H2").Value = "◆開発"

$ws.Range("A3").Value = "2026-01-08 06:31:27"
$ws.Range("B3").Value = "【急募】社内Webアプリの修正・再構築依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5467384"
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = "◇アプリ"

$ws.Range("A4").Value = "2026-01-08 06:31:27"
$ws.Range("B4").Value = "進行管理およびチームディレクションを担当"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "~ 5,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = "◇管理"

$ws.Range("A5").Value = "2026-01-08 06:31:27"
$ws.Range("B5").Value = "【急募】cloudflare導入の経験者を探しています!"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5467334"
$ws.Range("G5").Value = 18
$ws.Range("H5").ClearContents()

$ws.Range("A6").Value = "2026-01-08 06:31:27"
$ws.Range("B6").Value = "電気点火装置の回路図作成依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5466994"
$ws.Range("G6").Value = 13
$ws.Range("H6").ClearContents()

# --- Re-add hyperlinks for the URL column on the surviving rows ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5467295")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5467384")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5418064")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5467334")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5466994")

# --- Column width adjustments ---
# Excel's ColumnWidth (character units) differs from the OOXML <col width=".."/>
# by a constant ~0.83 padding offset on this font, so subtract it to land on
# the exact stored width (30 / 12) that the target file expects.
$ws.Columns.Item(2).ColumnWidth = 29.17
$ws.Columns.Item(8).ColumnWidth = 11.17
